# 21 de agosto de 2023 - Lap HP
# Se actualiza nuevamente el repositorio del curso de Fisica 1 Plan Cuatrimestral
#
# Updates physical-attendance / participation entries on "Concentrado" and the
# corresponding manually-copied grade column on "Parciales".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja "Concentrado": captura de asistencias / actividades faltantes.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Concentrado")

$ws.Range("P10").Value = 8
$ws.Range("P15").Value = 9
$ws.Range("Q17").Value = 9

$ws.Range("I25").Value = 1
$ws.Range("O25").Value = 9.2
$ws.Range("P25").Value = 9.2

$ws.Range("P33").Value = 9
$ws.Range("Q36").Value = 9

$ws.Range("F41").Value = 6
$ws.Range("G41").Value = 6

$ws.Range("F44").Value = 5
$ws.Range("G44").Value = 5
$ws.Range("J44").Value = 3
$ws.Range("M44").Value = 4

# View: el usuario cambio el zoom y la posicion de la ventana/seleccion antes
# de guardar.
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("R17").Select()

# Conditional formatting: las reglas de R y S se fusionan en un solo rango
# R3:S47 (antes eran dos reglas independientes).
$rRule = $ws.Range("R3:R47").FormatConditions.Item(1)
$sRule = $ws.Range("S3:S47").FormatConditions.Item(1)
$rRule.Delete()
$sRule.ModifyAppliesToRange($ws.Range("R3:S47"))

# ---------------------------------------------------------------------------
# Hoja "Parciales": se actualiza la calificacion del primer parcial (columna L,
# capturada manualmente) para que refleje los nuevos totales de "Concentrado".
# ---------------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("Parciales")

$wsP.Range("L3").Value = 7.7
$wsP.Range("L9").Value = 8.8
$wsP.Range("L14").Value = 9
$wsP.Range("L16").Value = 9
$wsP.Range("L17").Value = 9.3
$wsP.Range("L18").Value = 8.1
$wsP.Range("L20").Value = 7
$wsP.Range("L23").Value = 9.2
$wsP.Range("L29").Value = 8.4
$wsP.Range("L31").Value = 8.2
$wsP.Range("L34").Value = 8.5

$wsP.Activate()
$winP = $excel.ActiveWindow
$winP.Zoom = 120
$winP.ScrollRow = 7
$winP.ScrollColumn = 2
$wsP.Range("M9").Select()

$wb.Save()
